# Update the "Email" / "Password" values in the last data row (row 2) of the
# active worksheet with the latest combined test-case values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "sovereigncs032720_2015@test.com"
$ws.Range("L2").Value = "sovereigncs032720_2015"
